$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.611.85'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '2.008.73'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.632'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.90'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.55%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.390'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.75%  '

$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +16.20%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.880'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.13%  '

$ws.Range("D16").Value = '2.329.86'
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.86%  '

$ws.Range("D18").Value = '2.027.07'
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("D19").Value = '36.633.71'
$ws.Range("E19").Value = '  +1.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").Value = '0.0₃0873'
$ws.Range("E21").Value = '  +1.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.139'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +20.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.120'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0619'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.91%  '

$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("E39").Value = '  +0.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +21.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1000'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.85%  '

$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0214'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '93.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.03%  '

$ws.Range("D49").Value = '1.356.32'
$ws.Range("E49").Value = '  -5.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.31%  '

$ws.Range("D51").Value = '2.220.04'
$ws.Range("E51").Value = '  +1.04%  '
